$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "65.546.92"
$ws.Range("E2").Value = "  +2.02%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.644.89"
$ws.Range("E3").Value = "  +0.58%  "

# Row 4
$ws.Range("E4").Value = "  +0.13%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "603.84"
$ws.Range("E5").Value = "  +1.11%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "156.79"
$ws.Range("E6").Value = "  +2.99%  "

# Row 7
$ws.Range("E7").Value = "  +0.04%  "

# Row 8
$ws.Range("E8").Value = "  -0.59%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.642.83"
$ws.Range("E9").Value = "  +0.55%  "

# Row 10
$ws.Range("E10").Value = "  +7.71%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.402"
$ws.Range("E11").Value = "  +1.50%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "5.84"
$ws.Range("E12").Value = "  +0.08%  "

# Row 13
$ws.Range("E13").Value = "  +1.46%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "29.70"
$ws.Range("E14").Value = "  +5.74%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.0000194"
$ws.Range("E15").Value = "  +12.42%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.121.12"
$ws.Range("E16").Value = "  +0.71%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.256.27"
$ws.Range("E17").Value = "  +1.96%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.641.13"
$ws.Range("E18").Value = "  +1.74%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "12.62"
$ws.Range("E19").Value = "  +2.43%  "

# Row 20
$ws.Range("E20").Value = "  +1.49%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "356.86"
$ws.Range("E21").Value = "  +1.91%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.43"
$ws.Range("E22").Value = "  +4.33%  "

# Row 23
$ws.Range("E23").Value = "  +0.07%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "69.32"
$ws.Range("E24").Value = "  +2.42%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.71"
$ws.Range("E25").Value = "  +0.92%  "

# Row 26
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000105"
$ws.Range("E26").Value = "  +15.52%  "

# Row 27
$ws.Range("B27").Value = "InternetComputer(DFINITY)"
$ws.Range("C27").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.39"
$ws.Range("E27").Value = "  +1.47%  "

# Row 28
$ws.Range("E28").Value = "  -3.36%  "

# Row 29
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.165"
$ws.Range("E29").Value = "  +1.08%  "

# Row 30
$ws.Range("B30").Value = "Aptos"
$ws.Range("C30").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.13"
$ws.Range("E30").Value = "  -2.48%  "

# Row 31
$ws.Range("B31").Value = "Binance-PegBSC-USD"
$ws.Range("C31").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.00"
$ws.Range("E31").Value = "  +0.18%  "

# Row 32
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.17"
$ws.Range("E32").Value = "  +4.77%  "

# Row 33
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "529.24"
$ws.Range("E33").Value = "  -3.79%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.78"
$ws.Range("E34").Value = "  -3.16%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.55"
$ws.Range("E35").Value = "  +1.00%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.33"
$ws.Range("E36").Value = "  +2.17%  "

# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.430"
$ws.Range("E37").Value = "  +1.46%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "20.60"
$ws.Range("E38").Value = "  +1.97%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "161.67"
$ws.Range("E39").Value = "  -2.22%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.98"
$ws.Range("E40").Value = "  -1.14%  "

# Row 41
$ws.Range("E41").Value = "  +0.15%  "

# Row 42
$ws.Range("E42").Value = "  -0.02%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.02"
$ws.Range("E43").Value = "  +4.88%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "164.29"
$ws.Range("E44").Value = "  -3.04%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "4.14"
$ws.Range("E45").Value = "  +0.59%  "

# Row 46
$ws.Range("E46").Value = "  +6.57%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0608"
$ws.Range("E47").Value = "  +3.32%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "22.83"
$ws.Range("E48").Value = "  -2.20%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.651"
$ws.Range("E49").Value = "  +1.59%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0261"
$ws.Range("E50").Value = "  +3.03%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0982"
$ws.Range("E51").Value = "  +0.86%  "
